# Scheduled-runner update: refresh cached market-board profit figures
# (currentAveragePrice / LevePrice* / LeveProfit* columns H:N) across the
# per-job Sheets, per the latest price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 783.93335
$ws.Range("I41").Value = 976.9
$ws.Range("J41").Value = 398
$ws.Range("K41").Value = 976.9
$ws.Range("L41").Value = 398
$ws.Range("M41").Value = -536.9
$ws.Range("N41").Value = -1278

$ws.Range("H82").Value = 999.6
$ws.Range("I82").Value = 999.6
$ws.Range("K82").Value = 2998.8
$ws.Range("M82").Value = -2592.8

$ws.Range("H85").Value = 999.6
$ws.Range("I85").Value = 999.6
$ws.Range("K85").Value = 2998.8
$ws.Range("M85").Value = -1594.8

$ws.Range("H111").Value = 16909.666
$ws.Range("I111").Value = 16909.666
$ws.Range("K111").Value = 50728.99800000001
$ws.Range("M111").Value = -47661.99800000001

$ws.Range("H129").Value = 2109.9285
$ws.Range("I129").Value = 805.125
$ws.Range("J129").Value = 3849.6667
$ws.Range("K129").Value = 2415.375
$ws.Range("L129").Value = 11549.0001
$ws.Range("M129").Value = 2584.625
$ws.Range("N129").Value = -21549.0001

$ws.Range("H131").Value = 22253.143
$ws.Range("I131").Value = 962
$ws.Range("K131").Value = 2886
$ws.Range("M131").Value = 2154

$ws.Range("H137").Value = 1299.6
$ws.Range("I137").Value = 1236.1666
$ws.Range("J137").Value = 1394.75
$ws.Range("K137").Value = 3708.4998
$ws.Range("L137").Value = 4184.25
$ws.Range("M137").Value = -1158.4998
$ws.Range("N137").Value = -9284.25

$ws.Range("H141").Value = 2230.1875
$ws.Range("I141").Value = 2320.4285
$ws.Range("K141").Value = 6961.2855
$ws.Range("M141").Value = -1781.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3996.8096
$ws.Range("I32").Value = 3499.1843
$ws.Range("K32").Value = 3499.1843
$ws.Range("M32").Value = -3212.1843

$ws.Range("H61").Value = 1850.7646
$ws.Range("I61").Value = 1491.4445
$ws.Range("J61").Value = 2255
$ws.Range("K61").Value = 1491.4445
$ws.Range("L61").Value = 2255
$ws.Range("M61").Value = -1279.4445
$ws.Range("N61").Value = -2679

$ws.Range("H88").Value = 634.6
$ws.Range("I88").Value = 517.3333
$ws.Range("J88").Value = 684.8570999999999
$ws.Range("K88").Value = 517.3333
$ws.Range("L88").Value = 684.8570999999999
$ws.Range("M88").Value = -111.3333
$ws.Range("N88").Value = -1496.8571

$ws.Range("H91").Value = 634.6
$ws.Range("I91").Value = 517.3333
$ws.Range("J91").Value = 684.8570999999999
$ws.Range("K91").Value = 517.3333
$ws.Range("L91").Value = 684.8570999999999
$ws.Range("M91").Value = 886.6667
$ws.Range("N91").Value = -3492.8571

$ws.Range("H132").Value = 1548.3636
$ws.Range("I132").Value = 1502.5625
$ws.Range("K132").Value = 4507.6875
$ws.Range("M132").Value = -1977.6875

$ws.Range("H136").Value = 1850.7646
$ws.Range("I136").Value = 1491.4445
$ws.Range("J136").Value = 2255
$ws.Range("K136").Value = 4474.333500000001
$ws.Range("L136").Value = 6765
$ws.Range("M136").Value = -1924.333500000001
$ws.Range("N136").Value = -11865

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 830.2143
$ws.Range("J80").Value = 1002.8889
$ws.Range("L80").Value = 1002.8889
$ws.Range("N80").Value = -2998.8889

$ws.Range("H83").Value = 830.2143
$ws.Range("J83").Value = 1002.8889
$ws.Range("L83").Value = 5014.444500000001
$ws.Range("N83").Value = -14998.4445

$ws.Range("H107").Value = 113184.445
$ws.Range("I107").Value = 144380.14
$ws.Range("K107").Value = 144380.14
$ws.Range("M107").Value = -142460.14

$ws.Range("H134").Value = 1650.6522
$ws.Range("I134").Value = 1450.1052
$ws.Range("K134").Value = 4350.3156
$ws.Range("M134").Value = -1815.3156

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1739.0555
$ws.Range("I105").Value = 1861.8462
$ws.Range("J105").Value = 1419.8
$ws.Range("K105").Value = 1861.8462
$ws.Range("L105").Value = 1419.8
$ws.Range("M105").Value = -114.8462
$ws.Range("N105").Value = -4913.8

$ws.Range("H107").Value = 1580.3077
$ws.Range("I107").Value = 1817.0834
$ws.Range("K107").Value = 1817.0834
$ws.Range("M107").Value = 102.9166

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 50
$ws.Range("I7").Value = 2
$ws.Range("J7").Value = 98
$ws.Range("K7").Value = 6
$ws.Range("L7").Value = 294
$ws.Range("M7").Value = 106
$ws.Range("N7").Value = -518

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 23782
$ws.Range("J34").Value = 23782
$ws.Range("L34").Value = 23782
$ws.Range("N34").Value = -24318

$ws.Range("H76").Value = 23782
$ws.Range("J76").Value = 23782
$ws.Range("L76").Value = 23782
$ws.Range("N76").Value = -24412

$ws.Range("H79").Value = 23782
$ws.Range("J79").Value = 23782
$ws.Range("L79").Value = 23782
$ws.Range("N79").Value = -25966

$ws.Range("H80").Value = 3169.5
$ws.Range("I80").Value = 2066.5
$ws.Range("K80").Value = 2066.5
$ws.Range("M80").Value = -1068.5

$ws.Range("H83").Value = 3169.5
$ws.Range("I83").Value = 2066.5
$ws.Range("K83").Value = 10332.5
$ws.Range("M83").Value = -5340.5

$ws.Range("H102").Value = 2609.8462
$ws.Range("J102").Value = 1704
$ws.Range("L102").Value = 1704
$ws.Range("N102").Value = -4948

$ws.Range("H122").Value = 2132
$ws.Range("I122").Value = 2132
$ws.Range("K122").Value = 6396
$ws.Range("M122").Value = -3946

$ws.Range("H126").Value = 1435.2
$ws.Range("I126").Value = 1066.3334
$ws.Range("K126").Value = 3199.0002
$ws.Range("M126").Value = -729.0001999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 1120166.1

$ws.Range("H46").Value = 18808.385
$ws.Range("I46").Value = 37085
$ws.Range("J46").Value = 3142.7144
$ws.Range("K46").Value = 37085
$ws.Range("L46").Value = 3142.7144
$ws.Range("M46").Value = -36897
$ws.Range("N46").Value = -3518.7144

$ws.Range("H55").Value = 544.0769
$ws.Range("I55").Value = 668.4286
$ws.Range("J55").Value = 399
$ws.Range("K55").Value = 668.4286
$ws.Range("L55").Value = 399
$ws.Range("M55").Value = -495.4286
$ws.Range("N55").Value = -745

$ws.Range("H61").Value = 91387.87
$ws.Range("I61").Value = 85207.164
$ws.Range("K61").Value = 85207.164
$ws.Range("M61").Value = -85005.164

$ws.Range("H68").Value = 3037
$ws.Range("J68").Value = 3703
$ws.Range("L68").Value = 3703
$ws.Range("N68").Value = -5201

$ws.Range("H71").Value = 3037
$ws.Range("J71").Value = 3703
$ws.Range("L71").Value = 18515
$ws.Range("N71").Value = -26003

$ws.Range("H76").Value = 15531.333
$ws.Range("J76").Value = 15531.333
$ws.Range("L76").Value = 15531.333
$ws.Range("N76").Value = -16207.333

$ws.Range("H79").Value = 15531.333
$ws.Range("J79").Value = 15531.333
$ws.Range("L79").Value = 15531.333
$ws.Range("N79").Value = -17871.333

$ws.Range("H93").Value = 11877.781
$ws.Range("I93").Value = 1537.9642
$ws.Range("K93").Value = 1537.9642
$ws.Range("M93").Value = -289.9641999999999

$ws.Range("H113").Value = 91387.87
$ws.Range("I113").Value = 85207.164
$ws.Range("K113").Value = 85207.164
$ws.Range("M113").Value = -83037.164

$ws.Range("H132").Value = 3567.5
$ws.Range("I132").Value = 3243.5264
$ws.Range("K132").Value = 9730.5792
$ws.Range("M132").Value = -7200.5792

$ws.Range("H136").Value = 4603.9287
$ws.Range("I136").Value = 3446.625
$ws.Range("J136").Value = 6147
$ws.Range("K136").Value = 10339.875
$ws.Range("L136").Value = 18441
$ws.Range("M136").Value = -7789.875
$ws.Range("N136").Value = -23541

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 16668201
$ws.Range("I107").Value = 906.6316
$ws.Range("J107").Value = 45457164
$ws.Range("K107").Value = 2719.8948
$ws.Range("L107").Value = 136371492
$ws.Range("M107").Value = -799.8948
$ws.Range("N107").Value = -136375332

$ws.Range("H136").Value = 2276.6765
$ws.Range("I136").Value = 1891.5614
$ws.Range("K136").Value = 5674.6842
$ws.Range("M136").Value = -3124.6842
